# Add one more user (email + password) as a new row at the bottom of the
# usuario_mock sheet, matching the existing table's look & feel:
#  - column A (email) gets the same "Hyperlink" cell style as the other
#    email cells and a mailto: hyperlink
#  - column B (password) is plain text, no hyperlink (matches row 7's
#    pattern where a password that doesn't repeat "Password@123" is left
#    unstyled/unlinked)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEmail = "isaac@gmail.com"
$newPassword = "Password123$"

# New row is right after the current last row (row 7 -> row 8).
$emailCell = $ws.Range("A8")
$passwordCell = $ws.Range("B8")

$emailCell.Value = $newEmail
$passwordCell.Value = $newPassword

# Turn the email into a mailto hyperlink, then re-apply the same visual
# style already used by the other hyperlinked e-mail cells (e.g. A2) so we
# don't end up with a second, near-duplicate "Hyperlink" style.
$ws.Hyperlinks.Add($emailCell, "mailto:" + $newEmail)
$emailCell.Style = $ws.Range("A2").Style
